# feat: add 2022-Q4 data
#
# The workbook tracks JD (京东集团) holdings per quarter. This adds a new
# "2022-Q4" sheet (copied from the "2021-Q4" sheet's layout so headers /
# formatting line up) with its own fund rows, inserts it right after the
# "总计" (totals) summary sheet, and records the matching summary row on
# "总计" itself.

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item("总计")
$q2021  = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the "2022-Q4" detail sheet ---------------------------------
# Copy "2021-Q4" (placed right after "总计") so the header row / column
# widths / cell styling match the other quarterly sheets, then overwrite
# the data with the 2022-Q4 numbers.
$q2021.Copy($null, $totals)
$q2022 = $wb.Worksheets.Item(2)
$q2022.Name = "2022-Q4"

# Row 2: 161620 / 融通核心价值混合（QDII）A
$q2022.Cells.Item(2, 1).Value = 0
$q2022.Cells.Item(2, 2).NumberFormat = "@"
$q2022.Cells.Item(2, 2).Value = "161620"
$q2022.Cells.Item(2, 3).NumberFormat = "@"
$q2022.Cells.Item(2, 3).Value = "融通核心价值混合（QDII）A"
$q2022.Cells.Item(2, 4).NumberFormat = "@"
$q2022.Cells.Item(2, 4).Value = "0.51"
$q2022.Cells.Item(2, 5).NumberFormat = "@"
$q2022.Cells.Item(2, 5).Value = "65.37"
$q2022.Cells.Item(2, 6).NumberFormat = "@"
$q2022.Cells.Item(2, 6).Value = "2.86"
$q2022.Cells.Item(2, 7).NumberFormat = "@"
$q2022.Cells.Item(2, 7).Value = "0.0146"
$q2022.Cells.Item(2, 8).Value = 9

# Row 3: 014127 / 融通核心价值混合（QDII）C
$q2022.Cells.Item(3, 1).Value = 1
$q2022.Cells.Item(3, 2).NumberFormat = "@"
$q2022.Cells.Item(3, 2).Value = "014127"
$q2022.Cells.Item(3, 3).NumberFormat = "@"
$q2022.Cells.Item(3, 3).Value = "融通核心价值混合（QDII）C"
$q2022.Cells.Item(3, 4).NumberFormat = "@"
$q2022.Cells.Item(3, 4).Value = "0.02"
$q2022.Cells.Item(3, 5).NumberFormat = "@"
$q2022.Cells.Item(3, 5).Value = "65.37"
$q2022.Cells.Item(3, 6).NumberFormat = "@"
$q2022.Cells.Item(3, 6).Value = "2.86"
$q2022.Cells.Item(3, 7).NumberFormat = "@"
$q2022.Cells.Item(3, 7).Value = "0.0006"
$q2022.Cells.Item(3, 8).Value = 9

# --- 2. Add the matching summary row on "总计" -----------------------------
# Insert a new row 2 (pushing the existing 2021-Q4 / 2020-Q4 rows down) and
# fill it in with the 2022-Q4 totals, then refresh the running index in
# column A of the row that used to be "1" (2020-Q4), which is now "2".
$totals.Rows.Item(2).Insert()

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q4"
$totals.Cells.Item(2, 3).Value = 2
$totals.Cells.Item(2, 4).Value = 0.02

$totals.Cells.Item(4, 1).Value = 2

# Keep the originally-active tab ("2020-Q4") selected, since Copy()/rename
# above shifted focus onto the newly added sheet.
$wb.Worksheets.Item("2020-Q4").Activate()
